# Add a new "Smelt Larva Survey" (SLS) record as row 15 of Sheet1, below
# the existing 14 survey rows (A1:J14 -> A1:J15), matching the commit
# "Adding Smelt Larva Survey".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting first -------------------------------------------------
# Pull per-column formatting from existing cells that already carry the
# style the new row needs, via PasteSpecial(xlPasteFormats), so we reuse
# the workbook's existing style entries instead of minting new ones.
$ws.Range("A7").Copy()
$ws.Range("A15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H2").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C7").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("F15").PasteSpecial(-4122)

$ws.Range("G2").Copy()
$ws.Range("G15").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("H15").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("I15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Values -------------------------------------------------------------
# Survey, Survey_link, Abbreviation, Agency, Start year, Season,
# Frequency, Data source 1, Data_source_name
$ws.Range("C15").Value = "SLS"
$ws.Range("F15").Value = "Jan-Mar"
$ws.Range("H15").Value = "https://portal.edirepository.org/nis/mapbrowse?packageid=edi.534.2"
$ws.Range("B15").Value = "https://iep.ca.gov/Science-Synthesis-Service/Monitoring-Programs/Smelt-Larva"
$ws.Range("A15").Value = "Smelt Larva Survey"
$ws.Range("D15").Value = "California Department of Fish and Wildlife"
$ws.Range("E15").Value = 2009
$ws.Range("G15").Value = "2X monthly"
$ws.Range("I15").Value = "EDI"

# Selection ends up below the new row, like the authored edit.
$ws.Range("A16").Select()
